$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns A & B hold dd.mm.yyyy strings ("12.04.2024" etc). If written with
# Value while the cell is still in "General" format, Excel auto-detects them
# as dates and stores a serial number instead of the literal text - so we
# briefly force Text format, assign the literal string, then restore General
# so the cell style matches the rest of the sheet.
function Set-TextValue($range, [string]$text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.NumberFormat = "General"
}

# ---- Row 188 : 12.04.2024 ----
Set-TextValue $ws.Range("A188:B188") "12.04.2024"
$ws.Range("C188").Value = 33634
$ws.Range("D188").Value = 13000
$ws.Range("E188").Value = 8400
$ws.Range("F188").Value = 76214
$ws.Range("I188").Value = 8000
$ws.Range("J188").Value = 460
$ws.Range("K188").Value = 117
$ws.Range("L188").Value = 4750
$ws.Range("M188").Value = "https://web.archive.org/web/20240412232348/https://www.aljazeera.com/news/longform/2023/10/9/israel-hamas-war-in-maps-and-charts-live-tracker"

# ---- Row 189 : 13.04.2024 (figures unchanged vs 12.04.2024) ----
Set-TextValue $ws.Range("A189") "13.04.2024"
Set-TextValue $ws.Range("B189") "12.04.2024"
$ws.Range("C189").Value = 33634
$ws.Range("D189").Value = 13000
$ws.Range("E189").Value = 8400
$ws.Range("F189").Value = 76214
$ws.Range("I189").Value = 8000
$ws.Range("J189").Value = 460
$ws.Range("K189").Value = 117
$ws.Range("L189").Value = 4750
$ws.Range("M189").Value = "https://web.archive.org/web/20240413105351/https://www.aljazeera.com/news/longform/2023/10/9/israel-hamas-war-in-maps-and-charts-live-tracker"

# ---- Row 190 : 14.04.2024 ----
Set-TextValue $ws.Range("A190:B190") "14.04.2024"
$ws.Range("C190").Value = 33729
$ws.Range("D190").Value = 13800
$ws.Range("E190").Value = 8400
$ws.Range("F190").Value = 76371
$ws.Range("I190").Value = 8000
$ws.Range("J190").Value = 465
$ws.Range("K190").Value = 118
$ws.Range("L190").Value = 4750
$ws.Range("M190").Value = "https://web.archive.org/web/20240414152645/https://www.aljazeera.com/news/longform/2023/10/9/israel-hamas-war-in-maps-and-charts-live-tracker"

# ---- View state: mimic the author having scrolled down to the newly
# added rows, with the window split just above row 186 and the last
# selected cell being the final value entered (M190). ----
$win = $excel.ActiveWindow
$win.SplitRow = 185
$win.SplitColumn = 0
$ws.Range("D1").Select()
$ws.Range("M190").Select()
